$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 952.5714
$ws.Range("I33").Value = 944.6667
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 944.6667
$ws.Range("L33").Value = 1000
$ws.Range("M33").Value = -715.6667
$ws.Range("N33").Value = -1458
$ws.Range("H39").Value = 700
$ws.Range("H43").Value = 3080782.8
$ws.Range("J43").Value = 4766.3335
$ws.Range("L43").Value = 4766.3335
$ws.Range("N43").Value = -4904.3335
$ws.Range("H92").Value = 879
$ws.Range("I92").Value = 677.55554
$ws.Range("K92").Value = 677.55554
$ws.Range("M92").Value = 570.44446
$ws.Range("H106").Value = 23812146
$ws.Range("I106").Value = 30305098
$ws.Range("K106").Value = 30305098
$ws.Range("M106").Value = -30304467
$ws.Range("H129").Value = 1423.7407
$ws.Range("J129").Value = 2296.75
$ws.Range("L129").Value = 6890.25
$ws.Range("N129").Value = -16890.25
$ws.Range("H135").Value = 3924.4666
$ws.Range("I135").Value = 1273.6
$ws.Range("J135").Value = 9226.200000000001
$ws.Range("K135").Value = 11462.4
$ws.Range("L135").Value = 83035.8
$ws.Range("M135").Value = -8927.4
$ws.Range("N135").Value = -88105.8
$ws.Range("H137").Value = 16354151
$ws.Range("I137").Value = 2000920
$ws.Range("K137").Value = 6002760
$ws.Range("M137").Value = -6000210
$ws.Range("H138").Value = 4357.0527
$ws.Range("I138").Value = 1748.5555
$ws.Range("K138").Value = 5245.666499999999
$ws.Range("M138").Value = -105.6664999999994

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1590356.9
$ws.Range("I2").Value = 2497397.2
$ws.Range("J2").Value = 3036
$ws.Range("K2").Value = 2497397.2
$ws.Range("L2").Value = 3036
$ws.Range("M2").Value = -2497284.2
$ws.Range("N2").Value = -3262
$ws.Range("H32").Value = 12731.105
$ws.Range("I32").Value = 13357.951
$ws.Range("K32").Value = 13357.951
$ws.Range("M32").Value = -13070.951
$ws.Range("H45").Value = 3682.12
$ws.Range("I45").Value = 3484.389
$ws.Range("K45").Value = 3484.389
$ws.Range("M45").Value = -3107.389
$ws.Range("H46").Value = 12858.728
$ws.Range("I46").Value = 3250
$ws.Range("J46").Value = 14994
$ws.Range("K46").Value = 3250
$ws.Range("L46").Value = 14994
$ws.Range("M46").Value = -2931
$ws.Range("N46").Value = -15632
$ws.Range("H60").Value = 50000
$ws.Range("I60").Value = 50000
$ws.Range("K60").Value = 50000
$ws.Range("M60").Value = -49267
$ws.Range("H116").Value = 1590356.9
$ws.Range("I116").Value = 2497397.2
$ws.Range("J116").Value = 3036
$ws.Range("K116").Value = 2497397.2
$ws.Range("L116").Value = 3036
$ws.Range("M116").Value = -2495103.2
$ws.Range("N116").Value = -7624

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1590356.9
$ws.Range("I3").Value = 2497397.2
$ws.Range("J3").Value = 3036
$ws.Range("K3").Value = 2497397.2
$ws.Range("L3").Value = 3036
$ws.Range("M3").Value = -2497283.2
$ws.Range("N3").Value = -3264
$ws.Range("H88").Value = 17229.334
$ws.Range("J88").Value = 21344
$ws.Range("L88").Value = 21344
$ws.Range("N88").Value = -22156
$ws.Range("H91").Value = 17229.334
$ws.Range("J91").Value = 21344
$ws.Range("L91").Value = 21344
$ws.Range("N91").Value = -24152
$ws.Range("H94").Value = 762206.25
$ws.Range("I94").Value = 1142352.2
$ws.Range("K94").Value = 1142352.2
$ws.Range("M94").Value = -1141901.2
$ws.Range("H105").Value = 55557450
$ws.Range("I105").Value = 100001890
$ws.Range("J105").Value = 1901.125
$ws.Range("K105").Value = 100001890
$ws.Range("L105").Value = 1901.125
$ws.Range("M105").Value = -100000143
$ws.Range("N105").Value = -5395.125
$ws.Range("H107").Value = 2394.4443
$ws.Range("I107").Value = 2795.7144
$ws.Range("K107").Value = 2795.7144
$ws.Range("M107").Value = -875.7143999999998
$ws.Range("H134").Value = 3814.889
$ws.Range("I134").Value = 1608.091
$ws.Range("K134").Value = 4824.272999999999
$ws.Range("M134").Value = -2289.272999999999
$ws.Range("H138").Value = 93944.75
$ws.Range("J138").Value = 93944.75
$ws.Range("L138").Value = 93944.75
$ws.Range("N138").Value = -104224.75

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 992.38464
$ws.Range("I22").Value = 675.5
$ws.Range("J22").Value = 1499.4
$ws.Range("K22").Value = 675.5
$ws.Range("L22").Value = 1499.4
$ws.Range("M22").Value = -325.5
$ws.Range("N22").Value = -2199.4
$ws.Range("H31").Value = 5012.969
$ws.Range("I31").Value = 822.4516
$ws.Range("K31").Value = 822.4516
$ws.Range("M31").Value = -527.4516
$ws.Range("H34").Value = 5012.969
$ws.Range("I34").Value = 822.4516
$ws.Range("K34").Value = 822.4516
$ws.Range("M34").Value = -620.4516
$ws.Range("H58").Value = 307388.8
$ws.Range("I58").Value = 911185.9399999999
$ws.Range("K58").Value = 911185.9399999999
$ws.Range("M58").Value = -910982.9399999999
$ws.Range("H99").Value = 12634.333
$ws.Range("I99").Value = 21937.666
$ws.Range("K99").Value = 21937.666
$ws.Range("M99").Value = -20439.666
$ws.Range("H126").Value = 12634.333
$ws.Range("I126").Value = 21937.666
$ws.Range("K126").Value = 65812.99800000001
$ws.Range("M126").Value = -63342.99800000001
$ws.Range("H132").Value = 15163131
$ws.Range("I132").Value = 15885051
$ws.Range("K132").Value = 47655153
$ws.Range("M132").Value = -47652623
$ws.Range("H134").Value = 2201
$ws.Range("I134").Value = 1719.9474
$ws.Range("K134").Value = 5159.8422
$ws.Range("M134").Value = -2624.8422
$ws.Range("H136").Value = 307388.8
$ws.Range("I136").Value = 911185.9399999999
$ws.Range("K136").Value = 2733557.82
$ws.Range("M136").Value = -2731007.82

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 11481.889
$ws.Range("I3").Value = 4762.4287
$ws.Range("K3").Value = 14287.2861
$ws.Range("M3").Value = -14175.2861
$ws.Range("H4").Value = 6788078
$ws.Range("I4").Value = 2600289.5
$ws.Range("K4").Value = 7800868.5
$ws.Range("M4").Value = -7800756.5
$ws.Range("H14").Value = 267.7143
$ws.Range("I14").Value = 267.7143
$ws.Range("K14").Value = 803.1428999999999
$ws.Range("M14").Value = -630.1428999999999
$ws.Range("H18").Value = 4754.222
$ws.Range("J18").Value = 6833.3335
$ws.Range("L18").Value = 20500.0005
$ws.Range("N18").Value = -20838.0005
$ws.Range("H132").Value = 5572.364
$ws.Range("J132").Value = 9614.333000000001
$ws.Range("L132").Value = 86528.997
$ws.Range("N132").Value = -91588.997
$ws.Range("H133").Value = 31990
$ws.Range("I133").Value = 22990
$ws.Range("K133").Value = 68970
$ws.Range("M133").Value = -63910

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2169463
$ws.Range("I70").Value = 2981149.2
$ws.Range("J70").Value = 4966.6665
$ws.Range("K70").Value = 2981149.2
$ws.Range("L70").Value = 4966.6665
$ws.Range("M70").Value = -2980879.2
$ws.Range("N70").Value = -5506.6665
$ws.Range("H73").Value = 2169463
$ws.Range("I73").Value = 2981149.2
$ws.Range("J73").Value = 4966.6665
$ws.Range("K73").Value = 2981149.2
$ws.Range("L73").Value = 4966.6665
$ws.Range("M73").Value = -2980213.2
$ws.Range("N73").Value = -6838.6665
$ws.Range("H132").Value = 65957.30499999999
$ws.Range("I132").Value = 95505.32000000001
$ws.Range("J132").Value = 6861.273
$ws.Range("K132").Value = 286515.96
$ws.Range("L132").Value = 20583.819
$ws.Range("M132").Value = -283985.96
$ws.Range("N132").Value = -25643.819

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 48002010
$ws.Range("I16").Value = 70589736
$ws.Range("J16").Value = 3087.375
$ws.Range("K16").Value = 70589736
$ws.Range("L16").Value = 3087.375
$ws.Range("M16").Value = -70589566
$ws.Range("N16").Value = -3427.375
$ws.Range("H93").Value = 30666.334
$ws.Range("I93").Value = 29999
$ws.Range("J93").Value = 31000
$ws.Range("K93").Value = 29999
$ws.Range("L93").Value = 31000
$ws.Range("M93").Value = -28751
$ws.Range("N93").Value = -33496

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 300000
$ws.Range("J15").Value = 300000
$ws.Range("L15").Value = 300000
$ws.Range("N15").Value = -300576
$ws.Range("H100").Value = 1016221.56
$ws.Range("I100").Value = 1438050.1
$ws.Range("K100").Value = 2876100.2
$ws.Range("M100").Value = -2875559.2
$ws.Range("H122").Value = 4039.4424
$ws.Range("I122").Value = 3842.3333
$ws.Range("J122").Value = 4381.7896
$ws.Range("K122").Value = 11526.9999
$ws.Range("L122").Value = 13145.3688
$ws.Range("M122").Value = -9076.999899999999
$ws.Range("N122").Value = -18045.3688
$ws.Range("H126").Value = 2909.6191
$ws.Range("I126").Value = 2390
$ws.Range("J126").Value = 3948.8572
$ws.Range("K126").Value = 7170
$ws.Range("L126").Value = 11846.5716
$ws.Range("M126").Value = -4700
$ws.Range("N126").Value = -16786.5716
$ws.Range("H132").Value = 18522734
$ws.Range("I132").Value = 3268628.2
$ws.Range("J132").Value = 38470410
$ws.Range("K132").Value = 9805884.600000001
$ws.Range("L132").Value = 115411230
$ws.Range("M132").Value = -9803354.600000001
$ws.Range("N132").Value = -115416290
$ws.Range("H136").Value = 9051.35
$ws.Range("I136").Value = 2539.7646
$ws.Range("J136").Value = 11625.697
$ws.Range("K136").Value = 7619.293799999999
$ws.Range("L136").Value = 34877.091
$ws.Range("M136").Value = -5069.293799999999
$ws.Range("N136").Value = -39977.091
